$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(16, 8).Value = 5000
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = 5000
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 5000
$ws.Cells.Item(16, 13).Value = ""
$ws.Cells.Item(16, 14).Value = -5460

$ws.Cells.Item(113, 8).Value = 2007.125
$ws.Cells.Item(113, 9).Value = 1722.4286
$ws.Cells.Item(113, 11).Value = 1722.4286
$ws.Cells.Item(113, 13).Value = 1531.5714

$ws.Cells.Item(125, 8).Value = 781.7143
$ws.Cells.Item(125, 10).Value = 864
$ws.Cells.Item(125, 12).Value = 7776
$ws.Cells.Item(125, 14).Value = -12696

$ws.Cells.Item(132, 8).Value = 7252366.5
$ws.Cells.Item(132, 9).Value = 11496919
$ws.Cells.Item(132, 11).Value = 34490757
$ws.Cells.Item(132, 13).Value = -34488227

$ws.Cells.Item(138, 8).Value = 681133.75
$ws.Cells.Item(138, 9).Value = 1124.421
$ws.Cells.Item(138, 11).Value = 3373.263
$ws.Cells.Item(138, 13).Value = 1766.737

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(9, 8).Value = 7006
$ws.Cells.Item(9, 9).Value = 1000
$ws.Cells.Item(9, 11).Value = 1000
$ws.Cells.Item(9, 13).Value = -830

$ws.Cells.Item(20, 8).Value = 7006
$ws.Cells.Item(20, 9).Value = 1000
$ws.Cells.Item(20, 11).Value = 1000
$ws.Cells.Item(20, 13).Value = -730

$ws.Cells.Item(95, 8).Value = 15000
$ws.Cells.Item(95, 10).Value = 15000
$ws.Cells.Item(95, 12).Value = 15000
$ws.Cells.Item(95, 14).Value = -20492

$ws.Cells.Item(97, 8).Value = 994.3333
$ws.Cells.Item(97, 9).Value = 994.3333
$ws.Cells.Item(97, 11).Value = 994.3333
$ws.Cells.Item(97, 13).Value = -498.3333

$ws.Cells.Item(131, 8).Value = 39993.332
$ws.Cells.Item(131, 10).Value = 39993.332
$ws.Cells.Item(131, 12).Value = 39993.332
$ws.Cells.Item(131, 14).Value = -50073.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1803.174
$ws.Cells.Item(20, 9).Value = 1746.8572
$ws.Cells.Item(20, 11).Value = 1746.8572
$ws.Cells.Item(20, 13).Value = -1499.8572

$ws.Cells.Item(107, 8).Value = 1632.6957
$ws.Cells.Item(107, 9).Value = 1175.75
$ws.Cells.Item(107, 10).Value = 2131.182
$ws.Cells.Item(107, 11).Value = 1175.75
$ws.Cells.Item(107, 12).Value = 2131.182
$ws.Cells.Item(107, 13).Value = 744.25
$ws.Cells.Item(107, 14).Value = -5971.182

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 348.30768
$ws.Cells.Item(22, 9).Value = 348.9091
$ws.Cells.Item(22, 10).Value = 345
$ws.Cells.Item(22, 11).Value = 348.9091
$ws.Cells.Item(22, 12).Value = 345
$ws.Cells.Item(22, 13).Value = 1.090899999999976
$ws.Cells.Item(22, 14).Value = -1045

$ws.Cells.Item(31, 8).Value = 1651.5883
$ws.Cells.Item(31, 9).Value = 1505.1333
$ws.Cells.Item(31, 10).Value = 2750
$ws.Cells.Item(31, 11).Value = 1505.1333
$ws.Cells.Item(31, 12).Value = 2750
$ws.Cells.Item(31, 13).Value = -1210.1333
$ws.Cells.Item(31, 14).Value = -3340

$ws.Cells.Item(34, 8).Value = 1651.5883
$ws.Cells.Item(34, 9).Value = 1505.1333
$ws.Cells.Item(34, 10).Value = 2750
$ws.Cells.Item(34, 11).Value = 1505.1333
$ws.Cells.Item(34, 12).Value = 2750
$ws.Cells.Item(34, 13).Value = -1303.1333
$ws.Cells.Item(34, 14).Value = -3154

$ws.Cells.Item(132, 8).Value = 5915.885
$ws.Cells.Item(132, 9).Value = 7179.5293
$ws.Cells.Item(132, 10).Value = 3529
$ws.Cells.Item(132, 11).Value = 21538.5879
$ws.Cells.Item(132, 12).Value = 10587
$ws.Cells.Item(132, 13).Value = -19008.5879
$ws.Cells.Item(132, 14).Value = -15647

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(114, 8).Value = 722.35297
$ws.Cells.Item(114, 9).Value = 358.77777
$ws.Cells.Item(114, 10).Value = 1131.375
$ws.Cells.Item(114, 11).Value = 1076.33331
$ws.Cells.Item(114, 12).Value = 3394.125
$ws.Cells.Item(114, 13).Value = 2177.66669
$ws.Cells.Item(114, 14).Value = -9902.125

$ws.Cells.Item(122, 8).Value = 758.625
$ws.Cells.Item(122, 10).Value = 1200
$ws.Cells.Item(122, 12).Value = 10800
$ws.Cells.Item(122, 14).Value = -15700

$ws.Cells.Item(131, 8).Value = 931.76
$ws.Cells.Item(131, 9).Value = 419.625
$ws.Cells.Item(131, 10).Value = 976.29346
$ws.Cells.Item(131, 11).Value = 1258.875
$ws.Cells.Item(131, 12).Value = 2928.88038
$ws.Cells.Item(131, 13).Value = 3781.125
$ws.Cells.Item(131, 14).Value = -13008.88038

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 152.9
$ws.Cells.Item(2, 9).Value = 106.4
$ws.Cells.Item(2, 10).Value = 199.4
$ws.Cells.Item(2, 11).Value = 106.4
$ws.Cells.Item(2, 12).Value = 199.4
$ws.Cells.Item(2, 13).Value = 6.599999999999994
$ws.Cells.Item(2, 14).Value = -425.4

$ws.Cells.Item(3, 8).Value = 6833666.5
$ws.Cells.Item(3, 9).Value = 1000
$ws.Cells.Item(3, 11).Value = 1000
$ws.Cells.Item(3, 13).Value = -884

$ws.Cells.Item(10, 8).Value = 5000333.5
$ws.Cells.Item(10, 9).Value = 5000333.5
$ws.Cells.Item(10, 11).Value = 5000333.5
$ws.Cells.Item(10, 13).Value = -5000164.5

$ws.Cells.Item(11, 8).Value = 8344531
$ws.Cells.Item(11, 9).Value = 7660303
$ws.Cells.Item(11, 10).Value = 11252500
$ws.Cells.Item(11, 11).Value = 7660303
$ws.Cells.Item(11, 12).Value = 11252500
$ws.Cells.Item(11, 13).Value = -7660164
$ws.Cells.Item(11, 14).Value = -11252778

$ws.Cells.Item(24, 8).Value = 2502750
$ws.Cells.Item(24, 10).Value = 3666.6667
$ws.Cells.Item(24, 12).Value = 3666.6667
$ws.Cells.Item(24, 14).Value = -4012.6667

$ws.Cells.Item(80, 8).Value = 4015.9443
$ws.Cells.Item(80, 9).Value = 1817.1111
$ws.Cells.Item(80, 11).Value = 1817.1111
$ws.Cells.Item(80, 13).Value = -819.1111000000001

$ws.Cells.Item(83, 8).Value = 4015.9443
$ws.Cells.Item(83, 9).Value = 1817.1111
$ws.Cells.Item(83, 11).Value = 9085.5555
$ws.Cells.Item(83, 13).Value = -4093.5555

$ws.Cells.Item(86, 8).Value = 27997
$ws.Cells.Item(86, 10).Value = 27997
$ws.Cells.Item(86, 12).Value = 27997
$ws.Cells.Item(86, 14).Value = -30369

$ws.Cells.Item(89, 8).Value = 27997
$ws.Cells.Item(89, 10).Value = 27997
$ws.Cells.Item(89, 12).Value = 83991
$ws.Cells.Item(89, 14).Value = -95847

$ws.Cells.Item(95, 8).Value = 22745
$ws.Cells.Item(95, 10).Value = 22745
$ws.Cells.Item(95, 12).Value = 22745
$ws.Cells.Item(95, 14).Value = -28237

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1701.375
$ws.Cells.Item(7, 10).Value = 2765.5
$ws.Cells.Item(7, 12).Value = 2765.5
$ws.Cells.Item(7, 14).Value = -2989.5

$ws.Cells.Item(55, 8).Value = 439.33334
$ws.Cells.Item(55, 9).Value = 266
$ws.Cells.Item(55, 10).Value = 1046
$ws.Cells.Item(55, 11).Value = 266
$ws.Cells.Item(55, 12).Value = 1046
$ws.Cells.Item(55, 13).Value = -93
$ws.Cells.Item(55, 14).Value = -1392

$ws.Cells.Item(75, 8).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).Value = ""

$ws.Cells.Item(78, 8).Value = 0
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).Value = ""

$ws.Cells.Item(94, 8).Value = 15000
$ws.Cells.Item(94, 10).Value = 15000
$ws.Cells.Item(94, 12).Value = 15000
$ws.Cells.Item(94, 14).Value = -16352

$ws.Cells.Item(122, 8).Value = 22729826
$ws.Cells.Item(122, 9).Value = 41668564
$ws.Cells.Item(122, 10).Value = 3340.6
$ws.Cells.Item(122, 11).Value = 125005692
$ws.Cells.Item(122, 12).Value = 10021.8
$ws.Cells.Item(122, 13).Value = -125003242
$ws.Cells.Item(122, 14).Value = -14921.8

$ws.Cells.Item(126, 8).Value = 1701.375
$ws.Cells.Item(126, 10).Value = 2765.5
$ws.Cells.Item(126, 12).Value = 8296.5
$ws.Cells.Item(126, 14).Value = -13236.5

$ws.Cells.Item(132, 8).Value = 94690.09
$ws.Cells.Item(132, 9).Value = 3897.25
$ws.Cells.Item(132, 10).Value = 146571.72
$ws.Cells.Item(132, 11).Value = 11691.75
$ws.Cells.Item(132, 12).Value = 439715.16
$ws.Cells.Item(132, 13).Value = -9161.75
$ws.Cells.Item(132, 14).Value = -444775.16

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(48, 8).Value = 0
$ws.Cells.Item(48, 10).Value = 0
$ws.Cells.Item(48, 12).Value = 0
$ws.Cells.Item(48, 14).Value = ""

$ws.Cells.Item(97, 8).Value = 17857.334
$ws.Cells.Item(97, 10).Value = 17857.334
$ws.Cells.Item(97, 12).Value = 17857.334
$ws.Cells.Item(97, 14).Value = -19839.334
